$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 93-99), columns A:T (U only has the header note on row 1)
$data = @(
    @(45800, "Flowering",    "Large",  45, 55, 1.06, 1, "Yes", 2, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21),
    @(45800, "Nonflowering", "Medium", 45, 55, 1.06, 2, "Yes", 3, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21),
    @(45800, "Nonflowering", "Small",  45, 55, 1.06, 3, "Yes", 3, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21),
    @(45800, "Nonflowering", "Medium", 45, 55, 1.06, 3, "Yes", 3, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21),
    @(45800, "Nonflowering", "Medium", 45, 55, 1.06, 4, "Yes", 3, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21),
    @(45800, "Nonflowering", "Large",  45, 55, 1.06, 5, "Yes", 4, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21),
    @(45800, "Tree",         "Medium", 45, 55, 1.06, 6, "Yes", 1, "Neutral", 7, 0.64, 42, 29.98, 24, 0.9, 9.9, 23, 21)
)

$startRow = 93
$endRow = 99

# Copy the date-format style from the last existing date cell (A92) down
# the new rows so we reuse the workbook's existing style index instead of
# creating a duplicate custom number format.
$ws.Range("A92").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]            # A - Date
    $ws.Cells.Item($r, 2).Value = $row[1]            # B - Plant_Type
    $ws.Cells.Item($r, 3).Value = $row[2]            # C - Plant_Size
    $ws.Cells.Item($r, 4).Value = $row[3]            # D - Low
    $ws.Cells.Item($r, 5).Value = $row[4]            # E - High
    $ws.Cells.Item($r, 7).Value = $row[5]            # G - Rain
    $ws.Cells.Item($r, 8).Value = $row[6]            # H - Growth
    $ws.Cells.Item($r, 9).Value = $row[7]            # I - Pruned
    $ws.Cells.Item($r, 10).Value = $row[8]           # J - Quadrant
    $ws.Cells.Item($r, 11).Value = $row[9]           # K - Shade
    $ws.Cells.Item($r, 12).Value = $row[10]          # L - UV
    $ws.Cells.Item($r, 13).Value = $row[11]          # M - Humidity
    $ws.Cells.Item($r, 14).Value = $row[12]          # N - Dew_Point
    $ws.Cells.Item($r, 15).Value = $row[13]          # O - Pressure
    $ws.Cells.Item($r, 16).Value = $row[14]          # P - Wind_Gust
    $ws.Cells.Item($r, 17).Value = $row[15]          # Q - Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $row[16]          # R - Visibility
    $ws.Cells.Item($r, 19).Value = $row[17]          # S - AQI
    $ws.Cells.Item($r, 20).Value = $row[18]          # T - Pollen
}

# F - Temp_Diff: fill the formula across the whole new range in one go so
# Excel extends/creates the shared-formula group (matches how the existing
# F67:F92 block above it is built) rather than writing a literal formula
# into every single cell.
$ws.Range("F$startRow`:F$endRow").Formula = "=ABS(D$startRow-E$startRow)"

# Update the view: scroll down and select the newly added Wind_Gust column values
$ws.Application.ActiveWindow.ScrollRow = 91
$ws.Range("P93:P99").Select()
